# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 20 (pushing existing rows 20-152 down to 21-153)
# and populate it with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20; this shifts rows 20:152 down to 21:153
$ws.Range("A20:R20").EntireRow.Insert()

# Populate the newly inserted row 20 with this week's data
$ws.Cells.Item(20, 1).Value = 9
$ws.Cells.Item(20, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44462
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = 300000001
$ws.Cells.Item(20, 7).Value = "Rabanito"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 7900
$ws.Cells.Item(20, 11).Value = 3500
$ws.Cells.Item(20, 12).Value = 4000
$ws.Cells.Item(20, 13).Value = 3747
$ws.Cells.Item(20, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(20, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(20, 16).Value = 37
$ws.Cells.Item(20, 17).Value = 100
$ws.Cells.Item(20, 18).Value = "Hortaliza"
